# Add the team's season record (Wins / Losses / Ties) to the roster sheet.
# New columns AD:AF are appended after the existing AC ("Unnamed: 28") column,
# extending the used range from A1:AC48 to A1:AF48.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new labels, styled like the rest of the header row ---
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the existing header formatting (bold + centered + thin border) from the
# last header cell onto each new header cell so they match the other headers.
$ws.Range("AC1").Copy()
$ws.Range("AD1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("AC1").Copy()
$ws.Range("AE1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("AC1").Copy()
$ws.Range("AF1").PasteSpecial(-4122)   # xlPasteFormats

# --- Data rows (2-48): every player on this roster shares the team's 1998 record ---
for ($r = 2; $r -le 48; $r++) {
    $ws.Cells.Item($r, 30).Value = 65   # AD - Wins
    $ws.Cells.Item($r, 31).Value = 97   # AE - Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF - Ties
}
